$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Several "Price" values look numeric (e.g. "217.59"), but the source workbook stores
# them as plain text. Temporarily force a text format on each cell before assigning the
# new value so Excel does not silently convert it to a number, then clear the formatting
# again so the cell matches the original (unstyled) text cells.
$priceCells = @("D2","D3","D5","D8","D10","D11","D12","D13","D16","D17","D18","D19","D21","D22","D25","D29","D30","D34","D38","D39","D41","D43","D45","D46","D48","D49","D50","D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.739.82'
$ws.Range("D3").Value = '1.638.36'
$ws.Range("D5").Value = '217.59'
$ws.Range("D8").Value = '0.251'
$ws.Range("D10").Value = '19.08'
$ws.Range("D11").Value = '0.0844'
$ws.Range("D12").Value = '1.867.42'
$ws.Range("D13").Value = '1.644.27'
$ws.Range("D16").Value = '64.37'
$ws.Range("D17").Value = '26.737.53'
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("D19").Value = '211.56'
$ws.Range("D21").Value = '4.35'
$ws.Range("D22").Value = '2.35'
$ws.Range("D25").Value = '145.57'
$ws.Range("D29").Value = '15.59'
$ws.Range("D30").Value = '0.0505'
$ws.Range("D34").Value = '1.276.38'
$ws.Range("D38").Value = '0.531'
$ws.Range("D39").Value = '0.809'
$ws.Range("D41").Value = '0.802'
$ws.Range("D43").Value = '1.777.04'
$ws.Range("D45").Value = '60.78'
$ws.Range("D46").Value = '91.22'
$ws.Range("D48").Value = '0.0520'
$ws.Range("D49").Value = '7.55'
$ws.Range("D50").Value = '0.0962'
$ws.Range("D51").Value = '1.01'

foreach ($c in $priceCells) {
    $ws.Range($c).ClearFormats()
}

# --- Volume(1h) / Coin name / Link updates ---
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  -0.56%  '
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  +3.35%  '
$ws.Range("E24").Value = '  -2.65%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -0.76%  '
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("E37").Value = '  -1.59%  '
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("E42").Value = '  -1.85%  '
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("E44").Value = '  -3.53%  '
$ws.Range("E45").Value = '  +1.99%  '
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("E49").Value = '  -2.96%  '
$ws.Range("E50").Value = '  -0.45%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("E51").Value = '  +0.21%  '
